$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1573.1034
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1573.1034
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4719.3102
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5055.3102
$ws.Range("H92").Value = 191.47368
$ws.Range("I92").Value = 190.17647
$ws.Range("K92").Value = 190.17647
$ws.Range("M92").Value = 1057.82353
$ws.Range("H99").Value = 1171.5555
$ws.Range("I99").Value = 272.8
$ws.Range("J99").Value = 2295
$ws.Range("K99").Value = 818.4000000000001
$ws.Range("L99").Value = 6885
$ws.Range("M99").Value = 679.5999999999999
$ws.Range("N99").Value = -9881
$ws.Range("H101").Value = 2028.1666
$ws.Range("I101").Value = 1792.25
$ws.Range("K101").Value = 5376.75
$ws.Range("M101").Value = -3754.75
$ws.Range("H132").Value = 4429.5356
$ws.Range("I132").Value = 3536.05
$ws.Range("J132").Value = 6663.25
$ws.Range("K132").Value = 10608.15
$ws.Range("L132").Value = 19989.75
$ws.Range("M132").Value = -8078.150000000001
$ws.Range("N132").Value = -25049.75
$ws.Range("H137").Value = 316200.3
$ws.Range("I137").Value = 593073.4399999999
$ws.Range("J137").Value = 3213.3044
$ws.Range("K137").Value = 1779220.32
$ws.Range("L137").Value = 9639.913199999999
$ws.Range("M137").Value = -1776670.32
$ws.Range("N137").Value = -14739.9132

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 32204
$ws.Range("J80").Value = 37980
$ws.Range("L80").Value = 37980
$ws.Range("N80").Value = -39976
$ws.Range("H83").Value = 32204
$ws.Range("J83").Value = 37980
$ws.Range("L83").Value = 113940
$ws.Range("N83").Value = -123924
$ws.Range("H97").Value = 768.64703
$ws.Range("I97").Value = 717.8
$ws.Range("J97").Value = 1150
$ws.Range("K97").Value = 717.8
$ws.Range("L97").Value = 1150
$ws.Range("M97").Value = -221.8
$ws.Range("N97").Value = -2142
$ws.Range("H122").Value = 14096.875
$ws.Range("I122").Value = 23580.889
$ws.Range("K122").Value = 70742.667
$ws.Range("M122").Value = -68292.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 30624.666
$ws.Range("J35").Value = 30624.666
$ws.Range("L35").Value = 30624.666
$ws.Range("N35").Value = -31244.666
$ws.Range("H51").Value = 11750
$ws.Range("J51").Value = 11750
$ws.Range("L51").Value = 11750
$ws.Range("N51").Value = -12732
$ws.Range("H82").Value = 33838
$ws.Range("J82").Value = 34807.6
$ws.Range("L82").Value = 34807.6
$ws.Range("N82").Value = -35573.6
$ws.Range("H85").Value = 33838
$ws.Range("J85").Value = 34807.6
$ws.Range("L85").Value = 34807.6
$ws.Range("N85").Value = -37459.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2870.197
$ws.Range("I31").Value = 1321.225
$ws.Range("J31").Value = 5253.231
$ws.Range("K31").Value = 1321.225
$ws.Range("L31").Value = 5253.231
$ws.Range("M31").Value = -1026.225
$ws.Range("N31").Value = -5843.231
$ws.Range("H34").Value = 2870.197
$ws.Range("I34").Value = 1321.225
$ws.Range("J34").Value = 5253.231
$ws.Range("K34").Value = 1321.225
$ws.Range("L34").Value = 5253.231
$ws.Range("M34").Value = -1119.225
$ws.Range("N34").Value = -5657.231
$ws.Range("H51").Value = 9050.5
$ws.Range("J51").Value = 9050.5
$ws.Range("L51").Value = 9050.5
$ws.Range("N51").Value = -10522.5
$ws.Range("H60").Value = 25151.416
$ws.Range("J60").Value = 25151.416
$ws.Range("L60").Value = 25151.416
$ws.Range("N60").Value = -26173.416
$ws.Range("H61").Value = 9050.5
$ws.Range("J61").Value = 9050.5
$ws.Range("L61").Value = 9050.5
$ws.Range("N61").Value = -9746.5
$ws.Range("H76").Value = 2870
$ws.Range("I76").Value = 2870
$ws.Range("K76").Value = 2870
$ws.Range("M76").Value = -2555
$ws.Range("H79").Value = 2870
$ws.Range("I79").Value = 2870
$ws.Range("K79").Value = 2870
$ws.Range("M79").Value = -1778
$ws.Range("H105").Value = 1595
$ws.Range("I105").Value = 1090
$ws.Range("K105").Value = 1090
$ws.Range("M105").Value = 657
$ws.Range("H109").Value = 14329.167
$ws.Range("J109").Value = 14329.167
$ws.Range("L109").Value = 14329.167
$ws.Range("N109").Value = -16409.167
$ws.Range("H122").Value = 2930.923
$ws.Range("I122").Value = 1402
$ws.Range("J122").Value = 4241.4287
$ws.Range("K122").Value = 4206
$ws.Range("L122").Value = 12724.2861
$ws.Range("M122").Value = -1756
$ws.Range("N122").Value = -17624.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 5914.2856
$ws.Range("J101").Value = 5914.2856
$ws.Range("L101").Value = 17742.8568
$ws.Range("N101").Value = -22610.8568
$ws.Range("H131").Value = 875.3978
$ws.Range("J131").Value = 967.8553000000001
$ws.Range("L131").Value = 2903.5659
$ws.Range("N131").Value = -12983.5659
$ws.Range("H132").Value = 348037.22
$ws.Range("I132").Value = 1097546.4
$ws.Range("J132").Value = 2109.8845
$ws.Range("K132").Value = 9877917.6
$ws.Range("L132").Value = 18988.9605
$ws.Range("M132").Value = -9875387.6
$ws.Range("N132").Value = -24048.9605

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4813.2856
$ws.Range("J43").Value = 8739.714
$ws.Range("L43").Value = 8739.714
$ws.Range("N43").Value = -9041.714
$ws.Range("H57").Value = 13211.625
$ws.Range("J57").Value = 15930.5
$ws.Range("L57").Value = 15930.5
$ws.Range("N57").Value = -17570.5
$ws.Range("H122").Value = 2222.48
$ws.Range("I122").Value = 2150.1904
$ws.Range("J122").Value = 2602
$ws.Range("K122").Value = 6450.5712
$ws.Range("L122").Value = 7806
$ws.Range("M122").Value = -4000.5712
$ws.Range("N122").Value = -12706

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2844.8948
$ws.Range("I40").Value = 3109.8667
$ws.Range("J40").Value = 1851.25
$ws.Range("K40").Value = 3109.8667
$ws.Range("L40").Value = 1851.25
$ws.Range("M40").Value = -2973.8667
$ws.Range("N40").Value = -2123.25
$ws.Range("H93").Value = 1514.0741
$ws.Range("I93").Value = 1398.8572
$ws.Range("K93").Value = 1398.8572
$ws.Range("M93").Value = -150.8571999999999
$ws.Range("H136").Value = 2297.1714
$ws.Range("I136").Value = 2225.0344
$ws.Range("J136").Value = 2645.8333
$ws.Range("K136").Value = 6675.1032
$ws.Range("L136").Value = 7937.499899999999
$ws.Range("M136").Value = -4125.1032
$ws.Range("N136").Value = -13037.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 18688.5
$ws.Range("J109").Value = 18688.5
$ws.Range("L109").Value = 18688.5
$ws.Range("N109").Value = -21462.5
$ws.Range("H122").Value = 8269.941000000001
$ws.Range("I122").Value = 10862.909
$ws.Range("J122").Value = 3516.1667
$ws.Range("K122").Value = 32588.727
$ws.Range("L122").Value = 10548.5001
$ws.Range("M122").Value = -30138.727
$ws.Range("N122").Value = -15448.5001
